# The sheet originally held two contact rows (row2: 郭鑫惺, row3: 杜佳文).
# The commit replaces them with a list of mailboxes and appends a third
# row, growing the sheet from A1:B3 to A1:B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: QQ mailbox (replaces 郭鑫惺 / 2443760268gxx@gmail.com)
$ws.Range("A2").Value = "QQ邮箱"
$ws.Range("B2").Value = "1652675907@qq.com"

# Row 3: Google mailbox (replaces 杜佳文 / 1415006080@qq.com)
$ws.Range("A3").Value = "谷歌邮箱"
$ws.Range("B3").Value = "pitousanfadetuzi@gmail.com"

# Row 4: new NetEase mailbox row
$ws.Range("A4").Value = "网易邮箱"
$ws.Range("B4").Value = "13376003432@163.com"
